$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 325, pushing existing rows 325:345 down to 326:346
$ws.Rows.Item(325).Insert()

# Populate the newly inserted row 325 with the new record
$ws.Cells.Item(325, 1).Value = 3
$ws.Cells.Item(325, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(325, 3).Value = 'Coquimbo'
$ws.Cells.Item(325, 4).Value = (Get-Date -Year 2023 -Month 12 -Day 7 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(325, 5).Value = 5
$ws.Cells.Item(325, 6).Value = 100112026
$ws.Cells.Item(325, 7).Value = 'Haba'
$ws.Cells.Item(325, 8).Value = 'Sin especificar'
$ws.Cells.Item(325, 9).Value = 'Primera'
$ws.Cells.Item(325, 10).Value = 100
$ws.Cells.Item(325, 11).Value = 8500
$ws.Cells.Item(325, 12).Value = 9000
$ws.Cells.Item(325, 13).Value = 8750
$ws.Cells.Item(325, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(325, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(325, 16).Value = 350
$ws.Cells.Item(325, 17).Value = 25
$ws.Cells.Item(325, 18).Value = 'Hortaliza'
